$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: new match result added 08/08/2025 (U. Magdalena vs Pasto, draw 1-1)
# The date column holds plain text like "08/08/2025" elsewhere in the sheet,
# so force it with a leading apostrophe to stop Excel auto-converting it to
# a date serial, then restore the default "Normal" style so no stray
# number-format style sticks to the cell.
$ws.Range("A26").Value = "'08/08/2025"
$ws.Range("A26").Style = "Normal"

$ws.Range("B26").Value = "U. Magdalena"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = "Pasto"
$ws.Range("F26").Value = "D"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 1
$ws.Range("K26").Value = 0.73
$ws.Range("L26").Value = 0.68
$ws.Range("M26").Value = 12
$ws.Range("N26").Value = 8
$ws.Range("O26").Value = 7
$ws.Range("P26").Value = 2
